# The original paragraph held three runs:
#   "Label.text =”"  "Sreenija"  "”"
# It needs to collapse into a single run:
#   "Button.text  = “Bharath”"
$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "Label.text*Sreenija*") {
        # Build a Range spanning the whole paragraph (text only, not the
        # paragraph mark) and overwrite it; this merges/replaces every run
        # inside the paragraph with a single new run.
        $r = $d.Range($p.Range.Start, $p.Range.End)
        $r.Text = "Button.text  = “Bharath”"
    }
}
